# "Refonte des niveaux pour l'adaptation du procédural, 2 niveaux finis"
#
# The sheet gains a new column Q holding a small integer per enemy row
# (rows 3-10). Writing these values causes Excel/this runtime to extend
# the used range (dimension A2:P10 -> A2:Q10) and the per-row `spans`
# attribute automatically - no manual bookkeeping needed there.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("Q3").Value  = 2
$ws.Range("Q4").Value  = 1
$ws.Range("Q5").Value  = 2
$ws.Range("Q6").Value  = 3
$ws.Range("Q7").Value  = 1
$ws.Range("Q8").Value  = 4
$ws.Range("Q9").Value  = 4
$ws.Range("Q10").Value = 1

# The saved workbook ends with the selection parked on the last cell
# that was touched (Q10), matching the recorded <selection> element.
$ws.Range("Q10").Select() | Out-Null
